# Weekly data refresh: insert the new week's Albahaca price row at row 226
# (pushing all subsequent rows down by one) for the
# "Mercado Mayorista Lo Valledor de Santiago" wholesale market sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 226; Excel shifts row 226
# (and everything below it, through the old row 247) down to row 227..248.
$ws.Rows("226").Insert()

# Populate the newly inserted row 226 with the new week's record.
$ws.Range("A226").Value = 6
$ws.Range("B226").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C226").Value = "Metropolitana"
$ws.Range("D226").Value = 44461
$ws.Range("D226").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E226").Value = 13
$ws.Range("F226").Value = 100112052
$ws.Range("G226").Value = "Albahaca"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 150
$ws.Range("K226").Value = 5000
$ws.Range("L226").Value = 5500
$ws.Range("M226").Value = 5233
$ws.Range("N226").Value = "$/paquete"
$ws.Range("O226").Value = "Región de Arica y Parinacota"
$ws.Range("P226").Value = 5233
$ws.Range("Q226").Value = 1
$ws.Range("R226").Value = "Hortaliza"
